$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data refresh (Thu May  2 19:38:00 UTC 2024)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "59.166.45"
$ws.Range("E2").Value = "  +2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.989.66"
$ws.Range("E3").Value = "  +1.73%  "
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "561.80"
$ws.Range("E5").Value = "  +1.87%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.59"
$ws.Range("E6").Value = "  +6.34%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.979.59"
$ws.Range("E9").Value = "  +1.52%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.132"
$ws.Range("E10").Value = "  +3.82%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.10"
$ws.Range("E11").Value = "  +6.87%  "
$ws.Range("E12").Value = "  +2.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +3.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.61"
$ws.Range("E14").Value = "  +3.68%  "
$ws.Range("E15").Value = "  +2.23%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.490.51"
$ws.Range("E16").Value = "  +1.10%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.23"
$ws.Range("E17").Value = "  +7.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.991.98"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "59.146.34"
$ws.Range("E19").Value = "  +2.16%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "428.95"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.64"
$ws.Range("E21").Value = "  +5.00%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.717"
$ws.Range("E22").Value = "  +5.92%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.10"
$ws.Range("E23").Value = "  +2.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.29"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "80.75"
$ws.Range("E25").Value = "  +2.43%  "
$ws.Range("E26").Value = "  -0.04%  "
$ws.Range("E27").Value = "  -0.36%  "
$ws.Range("E28").Value = "  +9.33%  "
$ws.Range("E29").Value = "  +2.15%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.76"
$ws.Range("E30").Value = "  +4.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.69"
$ws.Range("E31").Value = "  +3.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.10"
$ws.Range("E32").Value = "  -0.13%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0986"
$ws.Range("E33").Value = "  -4.21%  "
$ws.Range("B34").Value = "Mantle"
$ws.Range("C34").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.993"
$ws.Range("E34").Value = "  +5.97%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.92"
$ws.Range("E35").Value = "  +6.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0761"
$ws.Range("E36").Value = "  +13.49%  "
$ws.Range("E37").Value = "  -0.67%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "48.96"
$ws.Range("E38").Value = "  +1.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.66"
$ws.Range("E39").Value = "  +4.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.70"
$ws.Range("E40").Value = "  +5.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "400.22"
$ws.Range("E41").Value = "  +6.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0350"
$ws.Range("E42").Value = "  +1.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.752.90"
$ws.Range("E43").Value = "  +4.46%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("E45").Value = "  +5.79%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "34.72"
$ws.Range("E47").Value = "  +21.91%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "121.31"
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.111"
$ws.Range("E49").Value = "  +1.78%  "
$ws.Range("E50").Value = "  +1.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "23.36"
$ws.Range("E51").Value = "  +0.84%  "
